# Update "想去人数" (number of interested attendees) counts in both the
# "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{Row = 5;  Old = 5195; New = 5199},
    @{Row = 6;  Old = 26;   New = 27},
    @{Row = 10; Old = 533;  New = 534},
    @{Row = 19; Old = 3159; New = 3161},
    @{Row = 27; Old = 132;  New = 133}
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Range("F" + $u.Row).Value = $u.New
}

$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @(
    @{Row = 6;  Old = 5195; New = 5199},
    @{Row = 7;  Old = 26;   New = 27},
    @{Row = 11; Old = 533;  New = 534},
    @{Row = 20; Old = 3159; New = 3161},
    @{Row = 28; Old = 132;  New = 133}
)
foreach ($u in $updates4) {
    $ws4.Range("F" + $u.Row).Value = $u.New
}
